$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 53.916668
$ws.Range("I6").Value2 = 31.545454
$ws.Range("J6").Value2 = 300
$ws.Range("K6").Value2 = 94.63636199999999
$ws.Range("L6").Value2 = 900
$ws.Range("M6").Value2 = 17.36363800000001
$ws.Range("N6").Value2 = -1124
$ws.Range("H17").Value2 = 3480.9
$ws.Range("J17").Value2 = 4481.2856
$ws.Range("L17").Value2 = 13443.8568
$ws.Range("N17").Value2 = -13779.8568
$ws.Range("H43").Value2 = 1300
$ws.Range("J43").Value2 = 0
$ws.Range("L43").Value2 = 0
$ws.Range("N43").ClearContents()
$ws.Range("H51").Value2 = 12500
$ws.Range("J51").Value2 = 12500
$ws.Range("L51").Value2 = 12500
$ws.Range("N51").Value2 = -13468
$ws.Range("H92").Value2 = 10418051
$ws.Range("I92").Value2 = 1402.8334
$ws.Range("J92").Value2 = 20834700
$ws.Range("K92").Value2 = 1402.8334
$ws.Range("L92").Value2 = 20834700
$ws.Range("M92").Value2 = -154.8334
$ws.Range("N92").Value2 = -20837196
$ws.Range("H116").Value2 = 4499
$ws.Range("J116").Value2 = 0
$ws.Range("L116").Value2 = 0
$ws.Range("N116").ClearContents()
$ws.Range("H127").Value2 = 15891.833
$ws.Range("I127").Value2 = 958.2857
$ws.Range("K127").Value2 = 2874.8571
$ws.Range("M127").Value2 = 2085.1429
$ws.Range("H129").Value2 = 145605.42
$ws.Range("I129").Value2 = 201958.1
$ws.Range("J129").Value2 = 4723.75
$ws.Range("K129").Value2 = 605874.3
$ws.Range("L129").Value2 = 14171.25
$ws.Range("M129").Value2 = -600874.3
$ws.Range("N129").Value2 = -24171.25
$ws.Range("H132").Value2 = 23548.227
$ws.Range("I132").Value2 = 851.2941
$ws.Range("K132").Value2 = 2553.8823
$ws.Range("M132").Value2 = -23.88229999999976
$ws.Range("H141").Value2 = 4341.1514
$ws.Range("I141").Value2 = 2343.6897
$ws.Range("K141").Value2 = 7031.0691
$ws.Range("M141").Value2 = -1851.0691

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2280.1304
$ws.Range("I2").Value2 = 2375.2354
$ws.Range("K2").Value2 = 2375.2354
$ws.Range("M2").Value2 = -2262.2354
$ws.Range("H32").Value2 = 13703273
$ws.Range("I32").Value2 = 15629160
$ws.Range("K32").Value2 = 15629160
$ws.Range("M32").Value2 = -15628873
$ws.Range("H45").Value2 = 3432.3333
$ws.Range("I45").Value2 = 3564
$ws.Range("K45").Value2 = 3564
$ws.Range("M45").Value2 = -3187
$ws.Range("H116").Value2 = 2280.1304
$ws.Range("I116").Value2 = 2375.2354
$ws.Range("K116").Value2 = 2375.2354
$ws.Range("M116").Value2 = -81.23540000000003
$ws.Range("H132").Value2 = 1268.1
$ws.Range("I132").Value2 = 1291.6171
$ws.Range("K132").Value2 = 3874.8513
$ws.Range("M132").Value2 = -1344.8513

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2280.1304
$ws.Range("I3").Value2 = 2375.2354
$ws.Range("K3").Value2 = 2375.2354
$ws.Range("M3").Value2 = -2261.2354
$ws.Range("H105").Value2 = 4792.385
$ws.Range("I105").Value2 = 2416.6667
$ws.Range("K105").Value2 = 2416.6667
$ws.Range("M105").Value2 = -669.6667000000002
$ws.Range("H107").Value2 = 1705.6875
$ws.Range("I107").Value2 = 1753.1538
$ws.Range("K107").Value2 = 1753.1538
$ws.Range("M107").Value2 = 166.8462
$ws.Range("H132").Value2 = 64748.977
$ws.Range("J132").Value2 = 64748.977
$ws.Range("L132").Value2 = 64748.977
$ws.Range("N132").Value2 = -74868.977

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value2 = 1750157.1
$ws.Range("I6").Value2 = 1750157.1
$ws.Range("K6").Value2 = 1750157.1
$ws.Range("M6").Value2 = -1750044.1
$ws.Range("H70").Value2 = 24999.5
$ws.Range("J70").Value2 = 24999.5
$ws.Range("L70").Value2 = 24999.5
$ws.Range("N70").Value2 = -25629.5
$ws.Range("H73").Value2 = 24999.5
$ws.Range("J73").Value2 = 24999.5
$ws.Range("L73").Value2 = 24999.5
$ws.Range("N73").Value2 = -27183.5
$ws.Range("H107").Value2 = 8706
$ws.Range("I107").Value2 = 1499.625
$ws.Range("M107").Value2 = 420.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value2 = 1222
$ws.Range("J7").Value2 = 1130.6
$ws.Range("L7").Value2 = 3391.8
$ws.Range("N7").Value2 = -3615.8
$ws.Range("H120").Value2 = 5558.3335
$ws.Range("I120").Value2 = 5558.3335
$ws.Range("K120").Value2 = 16675.0005
$ws.Range("M120").Value2 = -11837.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value2 = 2443200
$ws.Range("I11").Value2 = 1901142.9
$ws.Range("K11").Value2 = 1901142.9
$ws.Range("M11").Value2 = -1901003.9
$ws.Range("H15").Value2 = 25000
$ws.Range("J15").Value2 = 25000
$ws.Range("L15").Value2 = 25000
$ws.Range("N15").Value2 = -25576
$ws.Range("H80").Value2 = 30407268
$ws.Range("J80").Value2 = 41670132
$ws.Range("L80").Value2 = 41670132
$ws.Range("N80").Value2 = -41672128
$ws.Range("H81").Value2 = 25000
$ws.Range("J81").Value2 = 25000
$ws.Range("L81").Value2 = 25000
$ws.Range("N81").Value2 = -26996
$ws.Range("H83").Value2 = 30407268
$ws.Range("J83").Value2 = 41670132
$ws.Range("L83").Value2 = 208350660
$ws.Range("N83").Value2 = -208360644
$ws.Range("H84").Value2 = 25000
$ws.Range("J84").Value2 = 25000
$ws.Range("L84").Value2 = 75000
$ws.Range("N84").Value2 = -84984
$ws.Range("H102").Value2 = 1738.4
$ws.Range("I102").Value2 = 1150.921
$ws.Range("K102").Value2 = 1150.921
$ws.Range("M102").Value2 = 471.079
$ws.Range("H107").Value2 = 1155.3636
$ws.Range("I107").Value2 = 1388.625
$ws.Range("J107").Value2 = 533.3333
$ws.Range("K107").Value2 = 1388.625
$ws.Range("L107").Value2 = 533.3333
$ws.Range("M107").Value2 = 531.375
$ws.Range("N107").Value2 = -4373.3333
$ws.Range("H132").Value2 = 2779.7837
$ws.Range("I132").Value2 = 2736.6428
$ws.Range("K132").Value2 = 8209.928400000001
$ws.Range("M132").Value2 = -5679.928400000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2861
$ws.Range("J7").Value2 = 3519.3333
$ws.Range("L7").Value2 = 3519.3333
$ws.Range("N7").Value2 = -3743.3333
$ws.Range("H40").Value2 = 13727.111
$ws.Range("I40").Value2 = 16198.429
$ws.Range("K40").Value2 = 16198.429
$ws.Range("M40").Value2 = -16062.429
$ws.Range("H93").Value2 = 13892538
$ws.Range("I93").Value2 = 4045.5454
$ws.Range("K93").Value2 = 4045.5454
$ws.Range("M93").Value2 = -2797.5454
$ws.Range("H122").Value2 = 3452.1843
$ws.Range("I122").Value2 = 1877.65
$ws.Range("J122").Value2 = 5201.6665
$ws.Range("K122").Value2 = 5632.950000000001
$ws.Range("L122").Value2 = 15604.9995
$ws.Range("M122").Value2 = -3182.950000000001
$ws.Range("N122").Value2 = -20504.9995
$ws.Range("H126").Value2 = 2861
$ws.Range("J126").Value2 = 3519.3333
$ws.Range("L126").Value2 = 10557.9999
$ws.Range("N126").Value2 = -15497.9999
$ws.Range("H132").Value2 = 4291.7
$ws.Range("I132").Value2 = 3597.805
$ws.Range("J132").Value2 = 7452.778
$ws.Range("K132").Value2 = 10793.415
$ws.Range("L132").Value2 = 22358.334
$ws.Range("M132").Value2 = -8263.414999999999
$ws.Range("N132").Value2 = -27418.334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 8333.333000000001
$ws.Range("I14").Value2 = 5000
$ws.Range("J14").Value2 = 10000
$ws.Range("K14").Value2 = 5000
$ws.Range("L14").Value2 = 10000
$ws.Range("M14").Value2 = -4832
$ws.Range("N14").Value2 = -10336
$ws.Range("H70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("N73").ClearContents()
$ws.Range("H75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value2 = 0
$ws.Range("J78").Value2 = 0
$ws.Range("L78").Value2 = 0
$ws.Range("N78").ClearContents()
$ws.Range("H132").Value2 = 1311.5
$ws.Range("I132").Value2 = 1094.1818
$ws.Range("K132").Value2 = 3282.5454
$ws.Range("M132").Value2 = -752.5454
$ws.Range("H136").Value2 = 1711.3793
$ws.Range("I136").Value2 = 1199.5294
$ws.Range("J136").Value2 = 2436.5
$ws.Range("K136").Value2 = 3598.5882
$ws.Range("L136").Value2 = 7309.5
$ws.Range("M136").Value2 = -1048.5882
$ws.Range("N136").Value2 = -12409.5

